$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, pushing the existing "Total" row down to row 14.
$ws.Rows.Item(13).Insert()

# Copy formatting from row 12 (last data row) down into the new row 13
# so the new entry row looks/behaves like the other data rows.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new timesheet entry for row 13.
$ws.Range("A13").Value2 = 45274
$ws.Range("B13").Value2 = 0.583333333333333
$ws.Range("C13").Value2 = 0.833333333333333
$ws.Range("D13").Formula = "=(C13<B13)+C13-B13"
$ws.Range("E13").Value2 = 10
$ws.Range("F13").Formula = "=(D13*24)*E13"

# Re-apply the row's number formats, since entering the formulas above can
# make the engine infer a different (time-like) format for the result cells.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the Total row (now row 14) so its sums include the new row 13.
$ws.Range("D14").Formula = "=SUM(D2:D13)"
$ws.Range("F14").Formula = "=SUM(F2:F13)"

# Reflect the final cursor position left behind after entering the new data.
[void]$ws.Range("F15").Select()
